$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.352.74'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.801.83'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.36'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.576'
$ws.Range("E6").Value = '  +3.60%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '35.90'
$ws.Range("E8").Value = '  +9.41%  '
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0692'
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0964'
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("D12").Value = '2.062.97'
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.50'
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("D14").Value = '1.797.20'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.50'
$ws.Range("E16").Value = '  +4.99%  '
$ws.Range("D17").Value = '34.362.83'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.03'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.63'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '0.0₃0794'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("E21").Value = '  +2.17%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.17'
$ws.Range("E23").Value = '  +0.86%  '
$ws.Range("E24").Value = '  +3.33%  '
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.86'
$ws.Range("E26").Value = '  +7.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.86'
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("E28").Value = '  +2.70%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.04'
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("E31").Value = '  +1.19%  '
$ws.Range("E32").Value = '  +1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.82'
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").Value = '1.393.41'
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.672'
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.46'
$ws.Range("E37").Value = '  -4.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.23'
$ws.Range("E40").Value = '  +11.42%  '
$ws.Range("E41").Value = '  +2.68%  '
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '82.06'
$ws.Range("E43").Value = '  -2.62%  '
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.55'
$ws.Range("E45").Value = '  -3.94%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("E47").Value = '  -4.84%  '
$ws.Range("D48").Value = '1.963.25'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.73'
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  +0.74%  '
